# Updated cryptos list on Fri Sep  6 18:57:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell is written with a leading apostrophe so Excel stores it as literal
# text (matching the source inlineStr cells, many of which look numeric, e.g.
# "483.12" or "16.00" would otherwise be auto-coerced to a Number and lose
# trailing zeros). The Style reset clears the quote-prefix formatting bit that
# Excel applies to text-forced cells so no stray style/number-format diff is
# introduced relative to the original (unstyled) cells.

$ws.Range("D2").Value = "'53.527.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -5.23%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.223.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -6.58%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'483.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.99%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'125.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.20%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -5.58%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.231.24"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -6.61%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0917"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -7.13%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.42%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.32%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -3.53%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.622.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -6.49%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'21.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.84%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'53.438.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.32%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -3.88%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.228.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -6.35%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'Chainlink"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'9.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -5.21%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'Polkadot"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'3.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.30%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'297.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -3.07%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.09%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'63.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.50%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.363"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.46%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -4.23%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.55%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'169.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.17%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0680"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.21%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.04%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.17%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -4.30%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'17.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.26%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.73%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.841"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +5.42%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -6.04%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'35.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.42%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.364"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.21%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.27%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.83%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'122.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -6.50%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'4.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.09%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0875"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.64%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -5.89%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'229.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.15%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -3.16%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -3.67%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'16.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.51%  "
$ws.Range("E51").Style = "Normal"
